$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the C2 and C3 values from test-id-14 / test-id-15 to EMP-ID-25
$ws.Range("C2").Value = "EMP-ID-25"
$ws.Range("C3").Value = "EMP-ID-25"

# Update the active selection to C3 (was C4)
$ws.Range("C3").Select()
